$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap data rows 22 and 23 (match ids flip-flopped in source refresh) ---
$ws.Cells.Item(22, 2).Value = 6862618
$ws.Cells.Item(22, 3).Value = "Iraq League"
$ws.Cells.Item(22, 4).Value = 45112.45833333334
$ws.Cells.Item(22, 5).Value = "Naft AlWasat"
$ws.Cells.Item(22, 6).Value = "Newroz SC"
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 2
$ws.Cells.Item(22, 9).Value = "A"
$ws.Cells.Item(22, 10).Value = 3
$ws.Cells.Item(22, 11).Value = 2.4
$ws.Cells.Item(22, 12).Value = 2.75
$ws.Cells.Item(22, 13).Value = 2.05
$ws.Cells.Item(22, 14).Value = 2.625
$ws.Cells.Item(22, 15).Value = 4
$ws.Cells.Item(22, 16).Value = -0.25
$ws.Cells.Item(22, 17).Value = 1.775
$ws.Cells.Item(22, 18).Value = 2.025
$ws.Cells.Item(22, 19).Value = 2
$ws.Cells.Item(22, 20).Value = 1.95
$ws.Cells.Item(22, 21).Value = 1.85
$ws.Cells.Item(22, 22).Value = -1
$ws.Cells.Item(22, 23).Value = -1
$ws.Cells.Item(22, 24).Value = 3
$ws.Cells.Item(22, 25).Value = -1
$ws.Cells.Item(22, 26).Value = 1.025
$ws.Cells.Item(22, 27).Value = 0
$ws.Cells.Item(22, 28).Value = 0
$ws.Cells.Item(23, 2).Value = 6862617
$ws.Cells.Item(23, 3).Value = "Iraq League"
$ws.Cells.Item(23, 4).Value = 45112.45833333334
$ws.Cells.Item(23, 5).Value = "Al Naft SC"
$ws.Cells.Item(23, 6).Value = "Al Najaf"
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = "D"
$ws.Cells.Item(23, 10).Value = 3.6
$ws.Cells.Item(23, 11).Value = 2.8
$ws.Cells.Item(23, 12).Value = 2.1
$ws.Cells.Item(23, 13).Value = 1.95
$ws.Cells.Item(23, 14).Value = 2.8
$ws.Cells.Item(23, 15).Value = 4.1
$ws.Cells.Item(23, 16).Value = -0.5
$ws.Cells.Item(23, 17).Value = 2
$ws.Cells.Item(23, 18).Value = 1.8
$ws.Cells.Item(23, 19).Value = 2.25
$ws.Cells.Item(23, 20).Value = 2.025
$ws.Cells.Item(23, 21).Value = 1.775
$ws.Cells.Item(23, 22).Value = -1
$ws.Cells.Item(23, 23).Value = 1.8
$ws.Cells.Item(23, 24).Value = -1
$ws.Cells.Item(23, 25).Value = -1
$ws.Cells.Item(23, 26).Value = 0.8
$ws.Cells.Item(23, 27).Value = -1
$ws.Cells.Item(23, 28).Value = 0.7749999999999999

# --- Swap data rows 54 and 55 (match ids flip-flopped in source refresh) ---
$ws.Cells.Item(54, 2).Value = 7453404
$ws.Cells.Item(54, 3).Value = "Iraq League"
$ws.Cells.Item(54, 4).Value = 45240.375
$ws.Cells.Item(54, 5).Value = "Al Hudod"
$ws.Cells.Item(54, 6).Value = "Al Zawraa"
$ws.Cells.Item(54, 7).Value = 0
$ws.Cells.Item(54, 8).Value = 1
$ws.Cells.Item(54, 9).Value = "A"
$ws.Cells.Item(54, 10).Value = 3
$ws.Cells.Item(54, 11).Value = 2.9
$ws.Cells.Item(54, 12).Value = 2.3
$ws.Cells.Item(54, 13).Value = 3
$ws.Cells.Item(54, 14).Value = 2.9
$ws.Cells.Item(54, 15).Value = 2.3
$ws.Cells.Item(54, 16).Value = 0.25
$ws.Cells.Item(54, 17).Value = 1.75
$ws.Cells.Item(54, 18).Value = 2.05
$ws.Cells.Item(54, 19).Value = 2
$ws.Cells.Item(54, 20).Value = 1.9
$ws.Cells.Item(54, 21).Value = 1.9
$ws.Cells.Item(54, 22).Value = -1
$ws.Cells.Item(54, 23).Value = -1
$ws.Cells.Item(54, 24).Value = 1.3
$ws.Cells.Item(54, 25).Value = -1
$ws.Cells.Item(54, 26).Value = 1.05
$ws.Cells.Item(54, 27).Value = -1
$ws.Cells.Item(54, 28).Value = 0.8999999999999999
$ws.Cells.Item(55, 2).Value = 7453403
$ws.Cells.Item(55, 3).Value = "Iraq League"
$ws.Cells.Item(55, 4).Value = 45240.375
$ws.Cells.Item(55, 5).Value = "Al Qasim SC"
$ws.Cells.Item(55, 6).Value = "Naft AlBasra"
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(55, 9).Value = "D"
$ws.Cells.Item(55, 10).Value = 2.25
$ws.Cells.Item(55, 11).Value = 2.875
$ws.Cells.Item(55, 12).Value = 3.1
$ws.Cells.Item(55, 13).Value = 2.05
$ws.Cells.Item(55, 14).Value = 2.9
$ws.Cells.Item(55, 15).Value = 3.5
$ws.Cells.Item(55, 16).Value = -0.25
$ws.Cells.Item(55, 17).Value = 1.8
$ws.Cells.Item(55, 18).Value = 2
$ws.Cells.Item(55, 19).Value = 2
$ws.Cells.Item(55, 20).Value = 1.975
$ws.Cells.Item(55, 21).Value = 1.825
$ws.Cells.Item(55, 22).Value = -1
$ws.Cells.Item(55, 23).Value = 1.9
$ws.Cells.Item(55, 24).Value = -1
$ws.Cells.Item(55, 25).Value = -0.5
$ws.Cells.Item(55, 26).Value = 0.5
$ws.Cells.Item(55, 27).Value = -1
$ws.Cells.Item(55, 28).Value = 0.825

# --- Swap data rows 73 and 74 (match ids flip-flopped in source refresh) ---
$ws.Cells.Item(73, 2).Value = 7565124
$ws.Cells.Item(73, 3).Value = "Iraq League"
$ws.Cells.Item(73, 4).Value = 45268.35416666666
$ws.Cells.Item(73, 5).Value = "Al Naft SC"
$ws.Cells.Item(73, 6).Value = "Amanat Baghdad"
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 1
$ws.Cells.Item(73, 9).Value = "A"
$ws.Cells.Item(73, 10).Value = 1.5
$ws.Cells.Item(73, 11).Value = 3.6
$ws.Cells.Item(73, 12).Value = 6
$ws.Cells.Item(73, 13).Value = 1.727
$ws.Cells.Item(73, 14).Value = 3.5
$ws.Cells.Item(73, 15).Value = 4.2
$ws.Cells.Item(73, 16).Value = -0.75
$ws.Cells.Item(73, 17).Value = 2
$ws.Cells.Item(73, 18).Value = 1.8
$ws.Cells.Item(73, 19).Value = 2
$ws.Cells.Item(73, 20).Value = 1.875
$ws.Cells.Item(73, 21).Value = 1.925
$ws.Cells.Item(73, 22).Value = -1
$ws.Cells.Item(73, 23).Value = -1
$ws.Cells.Item(73, 24).Value = 3.2
$ws.Cells.Item(73, 25).Value = -1
$ws.Cells.Item(73, 26).Value = 0.8
$ws.Cells.Item(73, 27).Value = -1
$ws.Cells.Item(73, 28).Value = 0.925
$ws.Cells.Item(74, 2).Value = 7565123
$ws.Cells.Item(74, 3).Value = "Iraq League"
$ws.Cells.Item(74, 4).Value = 45268.35416666666
$ws.Cells.Item(74, 5).Value = "Al Karkh"
$ws.Cells.Item(74, 6).Value = "Karbalaa FC"
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 2
$ws.Cells.Item(74, 9).Value = "A"
$ws.Cells.Item(74, 10).Value = 1.727
$ws.Cells.Item(74, 11).Value = 3.4
$ws.Cells.Item(74, 12).Value = 4.2
$ws.Cells.Item(74, 13).Value = 1.833
$ws.Cells.Item(74, 14).Value = 3.3
$ws.Cells.Item(74, 15).Value = 3.75
$ws.Cells.Item(74, 16).Value = -0.5
$ws.Cells.Item(74, 17).Value = 1.875
$ws.Cells.Item(74, 18).Value = 1.925
$ws.Cells.Item(74, 19).Value = 2
$ws.Cells.Item(74, 20).Value = 1.925
$ws.Cells.Item(74, 21).Value = 1.875
$ws.Cells.Item(74, 22).Value = -1
$ws.Cells.Item(74, 23).Value = -1
$ws.Cells.Item(74, 24).Value = 2.75
$ws.Cells.Item(74, 25).Value = -1
$ws.Cells.Item(74, 26).Value = 0.925
$ws.Cells.Item(74, 27).Value = 0.925
$ws.Cells.Item(74, 28).Value = -1

# --- Swap data rows 91 and 92 (match ids flip-flopped in source refresh) ---
$ws.Cells.Item(91, 2).Value = 7618724
$ws.Cells.Item(91, 3).Value = "Iraq League"
$ws.Cells.Item(91, 4).Value = 45283.35416666666
$ws.Cells.Item(91, 5).Value = "Amanat Baghdad"
$ws.Cells.Item(91, 6).Value = "Karbalaa FC"
$ws.Cells.Item(91, 7).Value = 3
$ws.Cells.Item(91, 8).Value = 1
$ws.Cells.Item(91, 9).Value = "H"
$ws.Cells.Item(91, 10).Value = 2.875
$ws.Cells.Item(91, 11).Value = 2.75
$ws.Cells.Item(91, 12).Value = 2.5
$ws.Cells.Item(91, 13).Value = 2.875
$ws.Cells.Item(91, 14).Value = 2.75
$ws.Cells.Item(91, 15).Value = 2.5
$ws.Cells.Item(91, 16).Value = 0
$ws.Cells.Item(91, 17).Value = 2
$ws.Cells.Item(91, 18).Value = 1.8
$ws.Cells.Item(91, 19).Value = 2
$ws.Cells.Item(91, 20).Value = 2
$ws.Cells.Item(91, 21).Value = 1.8
$ws.Cells.Item(91, 22).Value = 1.875
$ws.Cells.Item(91, 23).Value = -1
$ws.Cells.Item(91, 24).Value = -1
$ws.Cells.Item(91, 25).Value = 1
$ws.Cells.Item(91, 26).Value = -1
$ws.Cells.Item(91, 27).Value = 1
$ws.Cells.Item(91, 28).Value = -1
$ws.Cells.Item(92, 2).Value = 7618726
$ws.Cells.Item(92, 3).Value = "Iraq League"
$ws.Cells.Item(92, 4).Value = 45283.35416666666
$ws.Cells.Item(92, 5).Value = "Al Hudod"
$ws.Cells.Item(92, 6).Value = "Duhok"
$ws.Cells.Item(92, 7).Value = 0
$ws.Cells.Item(92, 8).Value = 1
$ws.Cells.Item(92, 9).Value = "A"
$ws.Cells.Item(92, 10).Value = 2.375
$ws.Cells.Item(92, 11).Value = 2.8
$ws.Cells.Item(92, 12).Value = 3
$ws.Cells.Item(92, 13).Value = 2.375
$ws.Cells.Item(92, 14).Value = 2.8
$ws.Cells.Item(92, 15).Value = 3
$ws.Cells.Item(92, 16).Value = 0
$ws.Cells.Item(92, 17).Value = 1.75
$ws.Cells.Item(92, 18).Value = 2.05
$ws.Cells.Item(92, 19).Value = 1.75
$ws.Cells.Item(92, 20).Value = 1.8
$ws.Cells.Item(92, 21).Value = 2
$ws.Cells.Item(92, 22).Value = -1
$ws.Cells.Item(92, 23).Value = -1
$ws.Cells.Item(92, 24).Value = 2
$ws.Cells.Item(92, 25).Value = -1
$ws.Cells.Item(92, 26).Value = 1.05
$ws.Cells.Item(92, 27).Value = -1
$ws.Cells.Item(92, 28).Value = 1

# --- Swap data rows 108 and 109 (match ids flip-flopped in source refresh) ---
$ws.Cells.Item(108, 2).Value = 7811883
$ws.Cells.Item(108, 3).Value = "Iraq League"
$ws.Cells.Item(108, 4).Value = 45333.33333333334
$ws.Cells.Item(108, 5).Value = "Al Karkh"
$ws.Cells.Item(108, 6).Value = "Naft Maysan"
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 9).Value = "D"
$ws.Cells.Item(108, 10).Value = 4.2
$ws.Cells.Item(108, 11).Value = 3.4
$ws.Cells.Item(108, 12).Value = 1.727
$ws.Cells.Item(108, 13).Value = 3.6
$ws.Cells.Item(108, 14).Value = 3.3
$ws.Cells.Item(108, 15).Value = 1.909
$ws.Cells.Item(108, 16).Value = 0.5
$ws.Cells.Item(108, 17).Value = 1.825
$ws.Cells.Item(108, 18).Value = 1.975
$ws.Cells.Item(108, 19).Value = 1.75
$ws.Cells.Item(108, 20).Value = 1.725
$ws.Cells.Item(108, 21).Value = 1.975
$ws.Cells.Item(108, 22).Value = -1
$ws.Cells.Item(108, 23).Value = 2.3
$ws.Cells.Item(108, 24).Value = -1
$ws.Cells.Item(108, 25).Value = 0.825
$ws.Cells.Item(108, 26).Value = -1
$ws.Cells.Item(108, 27).Value = -1
$ws.Cells.Item(108, 28).Value = 0.9750000000000001
$ws.Cells.Item(109, 2).Value = 7811882
$ws.Cells.Item(109, 3).Value = "Iraq League"
$ws.Cells.Item(109, 4).Value = 45333.33333333334
$ws.Cells.Item(109, 5).Value = "Amanat Baghdad"
$ws.Cells.Item(109, 6).Value = "Naft AlWasat"
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = 1
$ws.Cells.Item(109, 9).Value = "D"
$ws.Cells.Item(109, 10).Value = 2.1
$ws.Cells.Item(109, 11).Value = 2.9
$ws.Cells.Item(109, 12).Value = 3.4
$ws.Cells.Item(109, 13).Value = 2.05
$ws.Cells.Item(109, 14).Value = 2.9
$ws.Cells.Item(109, 15).Value = 3.6
$ws.Cells.Item(109, 16).Value = -0.25
$ws.Cells.Item(109, 17).Value = 1.775
$ws.Cells.Item(109, 18).Value = 2.025
$ws.Cells.Item(109, 19).Value = 1.75
$ws.Cells.Item(109, 20).Value = 1.85
$ws.Cells.Item(109, 21).Value = 1.95
$ws.Cells.Item(109, 22).Value = -1
$ws.Cells.Item(109, 23).Value = 1.9
$ws.Cells.Item(109, 24).Value = -1
$ws.Cells.Item(109, 25).Value = -0.5
$ws.Cells.Item(109, 26).Value = 0.5125
$ws.Cells.Item(109, 27).Value = 0.425
$ws.Cells.Item(109, 28).Value = -0.5

# --- Append 7 new match rows (169-175) ---
$lastRow = 168
for ($i = 1; $i -le 7; $i++) {
    $src = $ws.Range("A" + $lastRow + ":AB" + $lastRow)
    $dst = $ws.Range("A" + ($lastRow + $i) + ":AB" + ($lastRow + $i))
    $src.Copy()
    $dst.PasteSpecial(-4122)
}

# Row 169
$ws.Cells.Item(169, 1).Value = 167
$ws.Cells.Item(169, 2).Value = 8122429
$ws.Cells.Item(169, 3).Value = "Iraq League"
$ws.Cells.Item(169, 4).Value = 45405.40625
$ws.Cells.Item(169, 5).Value = "Al Karkh"
$ws.Cells.Item(169, 6).Value = "Al Qasim SC"
$ws.Cells.Item(169, 7).Value = 1
$ws.Cells.Item(169, 8).Value = 1
$ws.Cells.Item(169, 9).Value = "D"
$ws.Cells.Item(169, 10).Value = 2.2
$ws.Cells.Item(169, 11).Value = 2.875
$ws.Cells.Item(169, 12).Value = 3.25
$ws.Cells.Item(169, 13).Value = 2.2
$ws.Cells.Item(169, 14).Value = 2.875
$ws.Cells.Item(169, 15).Value = 3.25
$ws.Cells.Item(169, 16).Value = -0.25
$ws.Cells.Item(169, 17).Value = 1.95
$ws.Cells.Item(169, 18).Value = 1.85
$ws.Cells.Item(169, 19).Value = 2
$ws.Cells.Item(169, 20).Value = 2
$ws.Cells.Item(169, 21).Value = 1.8
$ws.Cells.Item(169, 22).Value = -1
$ws.Cells.Item(169, 23).Value = 1.875
$ws.Cells.Item(169, 24).Value = -1
$ws.Cells.Item(169, 25).Value = -0.5
$ws.Cells.Item(169, 26).Value = 0.425
$ws.Cells.Item(169, 27).Value = 0
$ws.Cells.Item(169, 28).Value = 0

# Row 170
$ws.Cells.Item(170, 1).Value = 168
$ws.Cells.Item(170, 2).Value = 8131278
$ws.Cells.Item(170, 3).Value = "Iraq League"
$ws.Cells.Item(170, 4).Value = 45405.51041666666
$ws.Cells.Item(170, 5).Value = "Al Zawraa"
$ws.Cells.Item(170, 6).Value = "Naft AlBasra"
$ws.Cells.Item(170, 7).Value = 3
$ws.Cells.Item(170, 8).Value = 0
$ws.Cells.Item(170, 9).Value = "H"
$ws.Cells.Item(170, 10).Value = 1.615
$ws.Cells.Item(170, 11).Value = 3.2
$ws.Cells.Item(170, 12).Value = 5.5
$ws.Cells.Item(170, 13).Value = 1.533
$ws.Cells.Item(170, 14).Value = 3.4
$ws.Cells.Item(170, 15).Value = 5.75
$ws.Cells.Item(170, 16).Value = -1
$ws.Cells.Item(170, 17).Value = 2.025
$ws.Cells.Item(170, 18).Value = 1.775
$ws.Cells.Item(170, 19).Value = 1.75
$ws.Cells.Item(170, 20).Value = 1.725
$ws.Cells.Item(170, 21).Value = 1.975
$ws.Cells.Item(170, 22).Value = 0.5329999999999999
$ws.Cells.Item(170, 23).Value = -1
$ws.Cells.Item(170, 24).Value = -1
$ws.Cells.Item(170, 25).Value = 1.025
$ws.Cells.Item(170, 26).Value = -1
$ws.Cells.Item(170, 27).Value = 0.7250000000000001
$ws.Cells.Item(170, 28).Value = -1

# Row 171
$ws.Cells.Item(171, 1).Value = 169
$ws.Cells.Item(171, 2).Value = 8131279
$ws.Cells.Item(171, 3).Value = "Iraq League"
$ws.Cells.Item(171, 4).Value = 45405.60416666666
$ws.Cells.Item(171, 5).Value = "Al Najaf"
$ws.Cells.Item(171, 6).Value = "Amanat Baghdad"
$ws.Cells.Item(171, 7).Value = 1
$ws.Cells.Item(171, 8).Value = 1
$ws.Cells.Item(171, 9).Value = "D"
$ws.Cells.Item(171, 10).Value = 1.5
$ws.Cells.Item(171, 11).Value = 3.5
$ws.Cells.Item(171, 12).Value = 6
$ws.Cells.Item(171, 13).Value = 1.7
$ws.Cells.Item(171, 14).Value = 3.25
$ws.Cells.Item(171, 15).Value = 4.333
$ws.Cells.Item(171, 16).Value = -0.75
$ws.Cells.Item(171, 17).Value = 2
$ws.Cells.Item(171, 18).Value = 1.8
$ws.Cells.Item(171, 19).Value = 2
$ws.Cells.Item(171, 20).Value = 1.925
$ws.Cells.Item(171, 21).Value = 1.875
$ws.Cells.Item(171, 22).Value = -1
$ws.Cells.Item(171, 23).Value = 2.25
$ws.Cells.Item(171, 24).Value = -1
$ws.Cells.Item(171, 25).Value = -1
$ws.Cells.Item(171, 26).Value = 0.8
$ws.Cells.Item(171, 27).Value = 0
$ws.Cells.Item(171, 28).Value = 0

# Row 172
$ws.Cells.Item(172, 1).Value = 170
$ws.Cells.Item(172, 2).Value = 8140588
$ws.Cells.Item(172, 3).Value = "Iraq League"
$ws.Cells.Item(172, 4).Value = 45407.40625
$ws.Cells.Item(172, 5).Value = "Al Kahrabaa"
$ws.Cells.Item(172, 6).Value = "Zakho"
$ws.Cells.Item(172, 7).Value = 1
$ws.Cells.Item(172, 8).Value = 3
$ws.Cells.Item(172, 9).Value = "A"
$ws.Cells.Item(172, 10).Value = 3.4
$ws.Cells.Item(172, 11).Value = 2.75
$ws.Cells.Item(172, 12).Value = 2.2
$ws.Cells.Item(172, 13).Value = 2.6
$ws.Cells.Item(172, 14).Value = 2.6
$ws.Cells.Item(172, 15).Value = 3
$ws.Cells.Item(172, 16).Value = 0
$ws.Cells.Item(172, 17).Value = 1.75
$ws.Cells.Item(172, 18).Value = 2.05
$ws.Cells.Item(172, 19).Value = 1.75
$ws.Cells.Item(172, 20).Value = 1.825
$ws.Cells.Item(172, 21).Value = 1.975
$ws.Cells.Item(172, 22).Value = -1
$ws.Cells.Item(172, 23).Value = -1
$ws.Cells.Item(172, 24).Value = 2
$ws.Cells.Item(172, 25).Value = -1
$ws.Cells.Item(172, 26).Value = 1.05
$ws.Cells.Item(172, 27).Value = 0.825
$ws.Cells.Item(172, 28).Value = -1

# Row 173
$ws.Cells.Item(173, 1).Value = 171
$ws.Cells.Item(173, 2).Value = 8137758
$ws.Cells.Item(173, 3).Value = "Iraq League"
$ws.Cells.Item(173, 4).Value = 45407.5
$ws.Cells.Item(173, 5).Value = "Al Minaa"
$ws.Cells.Item(173, 6).Value = "Newroz SC"
$ws.Cells.Item(173, 7).Value = 1
$ws.Cells.Item(173, 8).Value = 4
$ws.Cells.Item(173, 9).Value = "A"
$ws.Cells.Item(173, 10).Value = 4.333
$ws.Cells.Item(173, 11).Value = 3.8
$ws.Cells.Item(173, 12).Value = 1.615
$ws.Cells.Item(173, 13).Value = 3.6
$ws.Cells.Item(173, 14).Value = 3.5
$ws.Cells.Item(173, 15).Value = 1.8
$ws.Cells.Item(173, 16).Value = 0.5
$ws.Cells.Item(173, 17).Value = 1.975
$ws.Cells.Item(173, 18).Value = 1.825
$ws.Cells.Item(173, 19).Value = 2.25
$ws.Cells.Item(173, 20).Value = 1.975
$ws.Cells.Item(173, 21).Value = 1.725
$ws.Cells.Item(173, 22).Value = -1
$ws.Cells.Item(173, 23).Value = -1
$ws.Cells.Item(173, 24).Value = 0.8
$ws.Cells.Item(173, 25).Value = -1
$ws.Cells.Item(173, 26).Value = 0.825
$ws.Cells.Item(173, 27).Value = 0.9750000000000001
$ws.Cells.Item(173, 28).Value = -1

# Row 174
$ws.Cells.Item(174, 1).Value = 172
$ws.Cells.Item(174, 2).Value = 8137759
$ws.Cells.Item(174, 3).Value = "Iraq League"
$ws.Cells.Item(174, 4).Value = 45407.5
$ws.Cells.Item(174, 5).Value = "Al Quwa Al Jawiya"
$ws.Cells.Item(174, 6).Value = "Karbalaa FC"
$ws.Cells.Item(174, 7).Value = 2
$ws.Cells.Item(174, 8).Value = 1
$ws.Cells.Item(174, 9).Value = "H"
$ws.Cells.Item(174, 10).Value = 1.571
$ws.Cells.Item(174, 11).Value = 3.6
$ws.Cells.Item(174, 12).Value = 5
$ws.Cells.Item(174, 13).Value = 1.25
$ws.Cells.Item(174, 14).Value = 4.5
$ws.Cells.Item(174, 15).Value = 10
$ws.Cells.Item(174, 16).Value = -1.75
$ws.Cells.Item(174, 17).Value = 2
$ws.Cells.Item(174, 18).Value = 1.8
$ws.Cells.Item(174, 19).Value = 2.75
$ws.Cells.Item(174, 20).Value = 1.925
$ws.Cells.Item(174, 21).Value = 1.875
$ws.Cells.Item(174, 22).Value = 0.25
$ws.Cells.Item(174, 23).Value = -1
$ws.Cells.Item(174, 24).Value = -1
$ws.Cells.Item(174, 25).Value = -1
$ws.Cells.Item(174, 26).Value = 0.8
$ws.Cells.Item(174, 27).Value = 0.4625
$ws.Cells.Item(174, 28).Value = -0.5

# Row 175
$ws.Cells.Item(175, 1).Value = 173
$ws.Cells.Item(175, 2).Value = 8136645
$ws.Cells.Item(175, 3).Value = "Iraq League"
$ws.Cells.Item(175, 4).Value = 45407.59375
$ws.Cells.Item(175, 5).Value = "Al Naft SC"
$ws.Cells.Item(175, 6).Value = "Al Shorta SC"
$ws.Cells.Item(175, 7).Value = 1
$ws.Cells.Item(175, 8).Value = 1
$ws.Cells.Item(175, 9).Value = "D"
$ws.Cells.Item(175, 10).Value = 3.6
$ws.Cells.Item(175, 11).Value = 3.2
$ws.Cells.Item(175, 12).Value = 1.909
$ws.Cells.Item(175, 13).Value = 4.5
$ws.Cells.Item(175, 14).Value = 3.3
$ws.Cells.Item(175, 15).Value = 1.7
$ws.Cells.Item(175, 16).Value = 0.75
$ws.Cells.Item(175, 17).Value = 1.825
$ws.Cells.Item(175, 18).Value = 1.975
$ws.Cells.Item(175, 19).Value = 2.25
$ws.Cells.Item(175, 20).Value = 1.875
$ws.Cells.Item(175, 21).Value = 1.925
$ws.Cells.Item(175, 22).Value = -1
$ws.Cells.Item(175, 23).Value = 2.3
$ws.Cells.Item(175, 24).Value = -1
$ws.Cells.Item(175, 25).Value = 0.825
$ws.Cells.Item(175, 26).Value = -1
$ws.Cells.Item(175, 27).Value = -0.5
$ws.Cells.Item(175, 28).Value = 0.4625
